$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update seed current density value from 0.006 to 0.005 A
# (dependent cells C11:C14, D10:D14, E10:E14, F10:F14 recalc via existing formulas)
$ws.Range("C10").Value = 0.005

# Update the active cell selection to F14, matching the saved state in the workbook
$ws.Range("F14").Select()
